{"js": "// The \"Programa\" section lists each topic prefixed with \"\u00a8 \" but all topics\n// were jammed together in a single run of text with no separators. This\n// edit inserts a manual line break (Word's vertical-tab / <w:br/>) before\n// every \"\u00a8 \" marker except the first one, so each topic renders on its own\n// line, for both the Portuguese and the italic English paragraphs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Vertical tab (U+000B) is how the Word object model represents a manual\n// line break (\"<w:br/>\") embedded inside a run's text.\nconst LINE_BREAK = \"\\u000B\";\n\nfunction withLineBreaks(text) {\n  // Split right before every \"\u00a8\" marker (keeping the marker), then re-join\n  // with a manual line break instead of nothing.\n  const parts = text.split(\"\u00a8\").filter((s) => s.length > 0);\n  return parts.map((s) => \"\u00a8\" + s).join(LINE_BREAK);\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  // Only touch paragraphs that still have two-or-more \"\u00a8 \" topics run\n  // together with no break between them yet (idempotent: a paragraph\n  // already containing a manual line break is left alone).\n  if (text.split(\"\u00a8\").filter((s) => s.length > 0).length < 2) continue;\n  if (text.indexOf(LINE_BREAK) !== -1) continue;\n\n  const newText = withLineBreaks(text);\n  if (newText === text) continue;\n\n  para.getRange().insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The \"Programa\" section lists each topic prefixed with \"\u00a8 \" but all topics\n# were jammed together in a single run of text with no separators. This\n# edit inserts a manual line break (vertical tab / <w:br/>) before every \"\u00a8 \"\n# marker except the first one, so each topic renders on its own line, for\n# both the Portuguese and the italic English paragraphs.\n\n$d = $word.ActiveDocument\n\n# Vertical tab (chr 11) is how the Word object model represents a manual\n# line break (\"<w:br/>\") embedded inside a run's text.\n$lineBreak = [char]11\n$marker = [char]0xA8\n$paraMark = [char]13\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text\n\n    # Only touch paragraphs that still have two-or-more \"\u00a8 \" topics run\n    # together with no break between them yet (idempotent: a paragraph\n    # already containing a manual line break is left alone).\n    $markerCount = ($text.ToCharArray() | Where-Object { $_ -eq $marker }).Count\n    if ($markerCount -lt 2) {\n        continue\n    }\n    if ($text.IndexOf($lineBreak) -ge 0) {\n        continue\n    }\n\n    # Paragraph.Range.Text carries the trailing paragraph mark (chr 13);\n    # strip it before splitting/rejoining so it doesn't get swallowed into\n    # the last segment and mint a spurious extra paragraph when written back.\n    $hadParaMark = $text.EndsWith($paraMark)\n    if ($hadParaMark) {\n        $text = $text.Substring(0, $text.Length - 1)\n    }\n\n    $parts = $text.Split($marker) | Where-Object { $_.Length -gt 0 }\n    $newText = ($parts | ForEach-Object { $marker + $_ }) -join $lineBreak\n\n    if ($newText -eq $text) {\n        continue\n    }\n\n    $p.Range.Text = $newText\n}\n"}
